# Automatische test-sync: 2025-06-23 18:43:50
#
# Adds the new "Verzoek om factuur" log entry (row 20) to the "Logs" sheet,
# extends the conditional formatting ranges to include the new row, and
# updates the "Dashboard" summary sheet so "Factuur / Administratie" moves
# ahead of "Offerte / Prijsaanvraag" with its incremented count.

$wb = $excel.ActiveWorkbook

# ---- Logs sheet: append the new row ----
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A20").Value = "Verzoek om factuur"
$logs.Range("B20").Value = "mailmind.test@zohomail.eu"
$logs.Range("C20").Value = "Kunt u mij een factuur sturen voor mijn laatste bestelling?"
$logs.Range("D20").Value = "Factuur / Administratie"
$logs.Range("F20").Value = "2025-06-23 18:43:47"
$logs.Range("G20").Value = "Nee"

# ---- Extend conditional formatting ranges to cover the new row ----
$catRules = $logs.Range("D2:D19").FormatConditions
$catRules.Item(1).ModifyAppliesToRange($logs.Range("D2:D20"))

$answeredRules = $logs.Range("G2:G19").FormatConditions
$answeredRules.Item(1).ModifyAppliesToRange($logs.Range("G2:G20"))

# ---- Dashboard sheet: swap rows 5/6 and bump the invoice count ----
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Range("A5").Value = "Factuur / Administratie"
$dashboard.Range("B5").Value = 3
$dashboard.Range("A6").Value = "Offerte / Prijsaanvraag"
$dashboard.Range("B6").Value = 2
